# Swap the contents of columns C (codeforiati:group-name) and D
# (codeforiati:group-code) for every row in the used range, including the
# header row, so that column C becomes the group-code and column D becomes
# the group-name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)

    $cValue = $cCell.Value()
    $dValue = $dCell.Value()

    $cCell.Value = $dValue
    $dCell.Value = $cValue
}
